$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.153.77"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "'2.516.49"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'589.70"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'178.05"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'0.143"
$ws.Range("E9").Value = "  +3.96%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'25.84"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "'67.983.96"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'2.850.82"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'2.479.21"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "'352.57"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "'4.12"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'70.86"
$ws.Range("E23").Value = "  +3.64%  "
$ws.Range("D24").Value = "'4.34"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'9.15"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "'2.642.24"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "'0.0₃0919"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "'509.49"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'7.89"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("D36").Value = "'164.25"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").Value = "'18.43"
$ws.Range("D38").Value = "'18.68"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  +5.26%  "
$ws.Range("D45").Value = "'147.24"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").Value = "'3.56"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'0.0₆0259"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("D49").Value = "'0.0745"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").Value = "'0.588"
$ws.Range("E51").Value = "  +0.99%  "
